# LicenseInfoResolver: Add original license source to ResolvedLicenseInfo
#
# The ResolvedLicense.toString() representation used by the report changed:
# the old "sources=[DECLARED]" fragment was replaced by a new
# "originalExpressions={DECLARED=[<license>]}" fragment (and it now comes
# after originalDeclaredLicenses instead of before it). This updates the
# three distinct "ResolvedLicense(...)" cell values on the
# "Gradle org.ossreviewtoolkit.gra" sheet (used by rows 12, 13/14 and 15)
# and highlights the touched cells with the same light-blue "resolved"
# background already used elsewhere in the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gradle org.ossreviewtoolkit.gra")

$lightBlue = 15128749  # RGB(173, 216, 230) == existing fillId=2 used across the sheet

function Update-ResolvedLicenseCell($cellRef, $newText) {
    $c = $ws.Range($cellRef)
    $c.Value = $newText
    $c.Interior.Color = $lightBlue
    $c.Font.Name = "Calibri"
    $c.Font.Size = 11
    $c.Font.Bold = $False
}

# Row 12 - Maven:junit:junit:4.12 - EPL-1.0
Update-ResolvedLicenseCell "C12" "ResolvedLicense(license=EPL-1.0, originalDeclaredLicenses=[Eclipse Public License 1.0], originalExpressions={DECLARED=[EPL-1.0]}, locations=[])"

# Row 13 - Maven:org.apache.commons:commons-lang3:3.5 - Apache-2.0
Update-ResolvedLicenseCell "C13" "ResolvedLicense(license=Apache-2.0, originalDeclaredLicenses=[Apache License, Version 2.0], originalExpressions={DECLARED=[Apache-2.0]}, locations=[])"

# Row 14 - Maven:org.apache.commons:commons-text:1.1 - Apache-2.0 (same text as row 13)
Update-ResolvedLicenseCell "C14" "ResolvedLicense(license=Apache-2.0, originalDeclaredLicenses=[Apache License, Version 2.0], originalExpressions={DECLARED=[Apache-2.0]}, locations=[])"

# Row 15 - Maven:org.hamcrest:hamcrest-core:1.3 - BSD-3-Clause
Update-ResolvedLicenseCell "C15" "ResolvedLicense(license=BSD-3-Clause, originalDeclaredLicenses=[New BSD License], originalExpressions={DECLARED=[BSD-3-Clause]}, locations=[])"

# Leave the cursor on the last-edited cell, as the author's session would have.
$ws.Activate()
$ws.Range("C15").Select()
